# Update cryptocurrency price (column D) and 1h volume change (column E)
# values on sheet1 to match the refreshed data pulled by the scheduled
# GitHub Actions job. Columns D/E are stored as text, so for values that
# look numeric we force text entry (leading apostrophe) and then clear the
# resulting 'Text' number format back to Normal so no stray styling is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.441.68"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "1.850.56"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'241.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").Value = "'0.6326"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.11%  "

$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "4.016.65"
$ws.Range("E8").Value = "  +113.38%  "

$ws.Range("D9").Value = "4.294.60"
$ws.Range("E9").Value = "  +98.67%  "

$ws.Range("D10").Value = "'0.07567"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").Value = "'0.2966"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").Value = "'24.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").Value = "'4.993"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.97%  "

$ws.Range("D15").Value = "'0.6850"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").Value = "'83.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("D17").Value = "'0.000009924"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.58%  "

$ws.Range("D18").Value = "'6.208"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "29.478.37"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "'231.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.59%  "

$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D23").Value = "'7.604"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "'155.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").Value = "'8.407"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("D29").Value = "4.219.59"
$ws.Range("E29").Value = "  +104.85%  "

$ws.Range("D30").Value = "'1.470"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").Value = "'0.05796"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.29%  "

$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").Value = "'4.132"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'4.022"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.35%  "

$ws.Range("D35").Value = "'1.857"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "

$ws.Range("E36").Value = "  -1.35%  "

$ws.Range("D37").Value = "'0.7170"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.74%  "

$ws.Range("E38").Value = "  -0.16%  "

$ws.Range("D39").Value = "1.252.00"
$ws.Range("E39").Value = "  +3.82%  "

$ws.Range("D40").Value = "'2.805"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("D41").Value = "'0.01804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.51%  "

$ws.Range("D42").Value = "'0.9025"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.93%  "

$ws.Range("D43").Value = "'6.103"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("D44").Value = "'0.9996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'101.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("D46").Value = "'67.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").Value = "'7.202"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("D48").Value = "'9.172"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").Value = "'1.684"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("D51").Value = "'0.1126"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "

